$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 35.333332
$ws.Range("I8").Value = 35.333332
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 105.999996
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 33.000004
$ws.Range("N8").Clear()

$ws.Range("H33").Value = 187.8
$ws.Range("I33").Value = 156.29411
$ws.Range("J33").Value = 366.33334
$ws.Range("K33").Value = 156.29411
$ws.Range("L33").Value = 366.33334
$ws.Range("M33").Value = 72.70589000000001
$ws.Range("N33").Value = -824.33334

$ws.Range("H41").Value = 894.625
$ws.Range("I41").Value = 680
$ws.Range("J41").Value = 992.1818
$ws.Range("K41").Value = 680
$ws.Range("L41").Value = 992.1818
$ws.Range("M41").Value = -240
$ws.Range("N41").Value = -1872.1818

$ws.Range("H62").Value = 2344.75
$ws.Range("J62").Value = 3003
$ws.Range("L62").Value = 3003
$ws.Range("N62").Value = -4251

$ws.Range("H65").Value = 2344.75
$ws.Range("J65").Value = 3003
$ws.Range("L65").Value = 15015
$ws.Range("N65").Value = -21255

$ws.Range("H70").Value = 3471.875
$ws.Range("J70").Value = 3464.1667
$ws.Range("L70").Value = 10392.5001
$ws.Range("N70").Value = -10932.5001

$ws.Range("H73").Value = 3471.875
$ws.Range("J73").Value = 3464.1667
$ws.Range("L73").Value = 10392.5001
$ws.Range("N73").Value = -12264.5001

$ws.Range("H95").Value = 87974
$ws.Range("J95").Value = 87974
$ws.Range("L95").Value = 87974
$ws.Range("N95").Value = -93466

$ws.Range("H98").Value = 314.89474
$ws.Range("I98").Value = 314.89474
$ws.Range("K98").Value = 314.89474
$ws.Range("M98").Value = 1183.10526

$ws.Range("H105").Value = 74975
$ws.Range("J105").Value = 74975
$ws.Range("L105").Value = 74975
$ws.Range("N105").Value = -81963

$ws.Range("H111").Value = 757.25
$ws.Range("I111").Value = 757.25
$ws.Range("K111").Value = 2271.75
$ws.Range("M111").Value = 795.25

$ws.Range("H122").Value = 314.89474
$ws.Range("I122").Value = 314.89474
$ws.Range("K122").Value = 944.6842200000001
$ws.Range("M122").Value = 1505.31578

$ws.Range("H128").Value = 150000
$ws.Range("J128").Value = 150000
$ws.Range("L128").Value = 150000
$ws.Range("N128").Value = -159960

$ws.Range("H129").Value = 1512.4
$ws.Range("I129").Value = 1512.4
$ws.Range("K129").Value = 4537.200000000001
$ws.Range("M129").Value = 462.7999999999993

$ws.Range("H135").Value = 1118.3103
$ws.Range("J135").Value = 1987.3334
$ws.Range("L135").Value = 17886.0006
$ws.Range("N135").Value = -22956.0006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2549
$ws.Range("I5").Value = 229.2
$ws.Range("K5").Value = 229.2
$ws.Range("M5").Value = -117.2

$ws.Range("H63").Value = 2149.5
$ws.Range("I63").Value = 2157.5264
$ws.Range("J63").Value = 1997
$ws.Range("K63").Value = 2157.5264
$ws.Range("L63").Value = 1997
$ws.Range("M63").Value = -1471.5264
$ws.Range("N63").Value = -3369

$ws.Range("H66").Value = 2149.5
$ws.Range("I66").Value = 2157.5264
$ws.Range("J66").Value = 1997
$ws.Range("K66").Value = 10787.632
$ws.Range("L66").Value = 9985
$ws.Range("M66").Value = -7355.632000000001
$ws.Range("N66").Value = -16849

$ws.Range("H74").Value = 5076.346
$ws.Range("I74").Value = 907.9286
$ws.Range("J74").Value = 22583.7
$ws.Range("K74").Value = 907.9286
$ws.Range("L74").Value = 22583.7
$ws.Range("M74").Value = -33.92859999999996
$ws.Range("N74").Value = -24331.7

$ws.Range("H77").Value = 5076.346
$ws.Range("I77").Value = 907.9286
$ws.Range("J77").Value = 22583.7
$ws.Range("K77").Value = 4539.643
$ws.Range("L77").Value = 112918.5
$ws.Range("M77").Value = -171.643
$ws.Range("N77").Value = -121654.5

$ws.Range("H122").Value = 3249.25
$ws.Range("J122").Value = 3499.5
$ws.Range("L122").Value = 10498.5
$ws.Range("N122").Value = -15398.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2549
$ws.Range("I4").Value = 229.2
$ws.Range("K4").Value = 229.2
$ws.Range("M4").Value = -114.2

$ws.Range("H80").Value = 1564.4
$ws.Range("I80").Value = 942
$ws.Range("K80").Value = 942
$ws.Range("M80").Value = 56

$ws.Range("H83").Value = 1564.4
$ws.Range("I83").Value = 942
$ws.Range("K83").Value = 4710
$ws.Range("M83").Value = 282

$ws.Range("H95").Value = 30208
$ws.Range("J95").Value = 30208
$ws.Range("L95").Value = 30208
$ws.Range("N95").Value = -35700

$ws.Range("H105").Value = 2425.9333
$ws.Range("J105").Value = 4080.6155
$ws.Range("L105").Value = 4080.6155
$ws.Range("N105").Value = -7574.6155

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").Clear()
$ws.Range("N107").Clear()

$ws.Range("H134").Value = 7140.7915
$ws.Range("I134").Value = 3017.2273
$ws.Range("K134").Value = 9051.6819
$ws.Range("M134").Value = -6516.6819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 13002.667
$ws.Range("I17").Value = 10504
$ws.Range("J17").Value = 18000
$ws.Range("K17").Value = 10504
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = -10330
$ws.Range("N17").Value = -18348

$ws.Range("H31").Value = 3060.75
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3060.75
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3060.75
$ws.Range("M31").Clear()
$ws.Range("N31").Value = -3650.75

$ws.Range("H34").Value = 3060.75
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3060.75
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3060.75
$ws.Range("M34").Clear()
$ws.Range("N34").Value = -3464.75

$ws.Range("H62").Value = 3633
$ws.Range("I62").Value = 3758.4
$ws.Range("J62").Value = 3006
$ws.Range("K62").Value = 3758.4
$ws.Range("L62").Value = 3006
$ws.Range("M62").Value = -3134.4
$ws.Range("N62").Value = -4254

$ws.Range("H65").Value = 3633
$ws.Range("I65").Value = 3758.4
$ws.Range("J65").Value = 3006
$ws.Range("K65").Value = 18792
$ws.Range("L65").Value = 15030
$ws.Range("M65").Value = -15672
$ws.Range("N65").Value = -21270

$ws.Range("H107").Value = 661.25
$ws.Range("I107").Value = 623.3077
$ws.Range("K107").Value = 623.3077
$ws.Range("M107").Value = 1296.6923

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2782.4285
$ws.Range("I12").Value = 3380
$ws.Range("J12").Value = 2334.25
$ws.Range("K12").Value = 10140
$ws.Range("L12").Value = 7002.75
$ws.Range("M12").Value = -9967
$ws.Range("N12").Value = -7348.75

$ws.Range("H68").Value = 6317.8184
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 6317.8184
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 18953.4552
$ws.Range("M68").Clear()
$ws.Range("N68").Value = -20575.4552

$ws.Range("H71").Value = 6317.8184
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 6317.8184
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 56860.3656
$ws.Range("M71").Clear()
$ws.Range("N71").Value = -64972.3656

$ws.Range("H80").Value = 2996

$ws.Range("H83").Value = 2996

$ws.Range("H92").Value = 420
$ws.Range("I92").Value = 340
$ws.Range("J92").Value = 900
$ws.Range("K92").Value = 1020
$ws.Range("L92").Value = 2700
$ws.Range("M92").Value = 228
$ws.Range("N92").Value = -5196

$ws.Range("H132").Value = 827.5833
$ws.Range("I132").Value = 580.3333
$ws.Range("K132").Value = 5222.9997
$ws.Range("M132").Value = -2692.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1517.15
$ws.Range("I80").Value = 1413.6
$ws.Range("J80").Value = 1620.7
$ws.Range("K80").Value = 1413.6
$ws.Range("L80").Value = 1620.7
$ws.Range("M80").Value = -415.5999999999999
$ws.Range("N80").Value = -3616.7

$ws.Range("H83").Value = 1517.15
$ws.Range("I83").Value = 1413.6
$ws.Range("J83").Value = 1620.7
$ws.Range("K83").Value = 7068
$ws.Range("L83").Value = 8103.5
$ws.Range("M83").Value = -2076
$ws.Range("N83").Value = -18087.5

$ws.Range("H105").Value = 23652
$ws.Range("J105").Value = 23652
$ws.Range("L105").Value = 23652
$ws.Range("N105").Value = -30640

$ws.Range("H107").Value = 4976.5415
$ws.Range("I107").Value = 8220.923000000001
$ws.Range("J107").Value = 1142.2727
$ws.Range("K107").Value = 8220.923000000001
$ws.Range("L107").Value = 1142.2727
$ws.Range("M107").Value = -6300.923000000001
$ws.Range("N107").Value = -4982.2727

$ws.Range("H126").Value = 3094.5386
$ws.Range("I126").Value = 2841.5
$ws.Range("J126").Value = 3499.4
$ws.Range("K126").Value = 8524.5
$ws.Range("L126").Value = 10498.2
$ws.Range("M126").Value = -6054.5
$ws.Range("N126").Value = -15438.2

$ws.Range("H132").Value = 9826.612999999999
$ws.Range("I132").Value = 11044.12
$ws.Range("K132").Value = 33132.36
$ws.Range("M132").Value = -30602.36

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 3139.2856
$ws.Range("I13").Value = 1842.3077
$ws.Range("K13").Value = 1842.3077
$ws.Range("M13").Value = -1702.3077

$ws.Range("H55").Value = 1802.7188
$ws.Range("I55").Value = 1739.4166
$ws.Range("J55").Value = 1840.7
$ws.Range("K55").Value = 1739.4166
$ws.Range("L55").Value = 1840.7
$ws.Range("M55").Value = -1566.4166
$ws.Range("N55").Value = -2186.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4790
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 5000
$ws.Range("N96").Value = -7746

$ws.Range("H113").Value = 967.75
$ws.Range("I113").Value = 1228.4286
$ws.Range("J113").Value = 602.8
$ws.Range("K113").Value = 3685.2858
$ws.Range("L113").Value = 1808.4
$ws.Range("M113").Value = -1515.2858
$ws.Range("N113").Value = -6148.4
